# Project Sample Project is saved. Author: admin. Type: SAVE.
# Rule "R40" (row 11 of the Rules decision table) is renamed to "1".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to keep the numeric-looking entry as
# text (same as a user typing '1 into the cell), matching the original
# "Rule" column which stores rule names as text.
$ws.Range("B11").Value = "'1"
